# Update the imputed values produced by the KNN algorithm for the
# terrestrial_mammals / combination_2_ABCDE / BCD / 10 / seed4 dataset.
# Only the specific cells listed below changed value; everything else
# (headers, formatting, other data) stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value  = -12.965
$ws.Range("B9").Value  = 5.317
$ws.Range("C12").Value = -11.171
$ws.Range("D15").Value = -8.463000000000001
$ws.Range("B18").Value = 5.282999999999999
$ws.Range("B20").Value = 7.242999999999999
$ws.Range("C26").Value = -13.131
$ws.Range("B27").Value = 5.526
$ws.Range("C27").Value = -13.664
$ws.Range("C29").Value = -12.335
$ws.Range("C37").Value = -13.351
$ws.Range("C38").Value = -13.738
$ws.Range("D38").Value = -7.473000000000001
$ws.Range("D44").Value = -7.568999999999998
$ws.Range("C51").Value = -12.267
$ws.Range("D51").Value = -7.542999999999999
$ws.Range("C55").Value = -13.752
$ws.Range("D57").Value = -8.130999999999998
$ws.Range("D63").Value = -7.337000000000001
$ws.Range("B69").Value = 5.667
$ws.Range("C69").Value = -11.17
$ws.Range("C70").Value = -12.338
$ws.Range("D70").Value = -7.798
$ws.Range("B76").Value = 6.723999999999999
$ws.Range("B82").Value = 5.343999999999999
$ws.Range("C83").Value = -13.551
$ws.Range("D99").Value = -7.512
$ws.Range("C102").Value = -13.419
